# Carbon journey cost added: new column BS "Carbon Kg per litre Fuel"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InputSchedule")

# Clone the formatting of column BR (header + 4 data rows) onto new column BS,
# then overwrite the values/header text. Using xlPasteAll (rather than
# xlPasteFormats) so the style index is reliably carried even when it is
# visually equivalent to the default style.
$ws.Range("BR1:BR5").Copy() | Out-Null
$ws.Range("BS1:BS5").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("BS1").Value = "Carbon Kg per litre Fuel"
$ws.Range("BS2").Value = 2.31
$ws.Range("BS3").Value = 2.31
$ws.Range("BS4").Value = 2.31
$ws.Range("BS5").Value = 2.31

# Extend the conditional formatting that already covers BQ2:BR4 and X2:AC4
# so that it also covers the new BS2:BS5 cells, matching the same rule
# (expression "X2<>X1", stop-if-true, same highlight fill).
$newCf = $ws.Range("BS2:BS5").FormatConditions.Add(2, 3, "=X2<>X1")
$newCf = $ws.Range("BS2:BS5").FormatConditions.Item($ws.Range("BS2:BS5").FormatConditions.Count)
$newCf.StopIfTrue = $true
$newCf.Interior.Color = 11854022
